$d = $word.ActiveDocument

# Change 1 (run-local): "What stands out for you about the " -> "Wat staan vir jou uit oor die "
$d.Content.Find.Execute("What stands out for you about the ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Wat staan vir jou uit oor die ", 2)

# Change 2 (run-local): "online safety" -> "aanlyn-veiligheid"
$d.Content.Find.Execute("online safety", $true, $false, $false, $false, $false,
                         $true, 1, $false, "aanlyn-veiligheid", 2)

# Change 3 (run-local): the lone space run before "lesson" -> "-"
# Restrict the search range to just before "lesson" so we do not touch any
# other single-space run elsewhere in the document.
$rng = $d.Content
$rng.Find.Execute("lesson", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$spaceRng = $d.Range($rng.Start - 1, $rng.Start)
$spaceRng.Find.Execute(" ", $true, $false, $false, $false, $false,
                        $true, 1, $false, "-", 2)

# Change 4 (run-local): "lesson" -> "les"
$d.Content.Find.Execute("lesson", $true, $false, $false, $false, $false,
                         $true, 1, $false, "les", 2)

# Change 5 (run-local): "Would you mind sharing with me what that was?" -> "Sal jy omgee om met my te deel wat dit was?"
$d.Content.Find.Execute("Would you mind sharing with me what that was?", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sal jy omgee om met my te deel wat dit was?", 2)
